$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("ELEMENTS")

# ---------------------------------------------------------------------------
# Block 1: Reducers -> final rows 3 (Conc.) and 4 (Ecc.), typed row by row.
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value = "EN 10253-2 - Reducer: Conc."
$ws.Cells.Item(3, 2).Value = "RED"
$ws.Cells.Item(3, 4).Value = "EN 10253-2"

$ws.Cells.Item(4, 1).Value = "EN 10253-2 - Reducer: Ecc."
$ws.Cells.Item(4, 2).Value = "RED"
$ws.Cells.Item(4, 4).Value = "EN 10253-2"

# ---------------------------------------------------------------------------
# Block 2: Flange weld collars -> final rows 5..9 (PN06, PN10, PN16, PN25,
# PN40). Column A entered first in PN40/PN06/PN10/PN16/PN25 order (written
# straight to each row's final position), then column B, then column D.
# ---------------------------------------------------------------------------
$ws.Cells.Item(9, 1).Value = "Flange weld collar: PN40"
$ws.Cells.Item(5, 1).Value = "Flange weld collar: PN06"
$ws.Cells.Item(6, 1).Value = "Flange weld collar: PN10"
$ws.Cells.Item(7, 1).Value = "Flange weld collar: PN16"
$ws.Cells.Item(8, 1).Value = "Flange weld collar: PN25"

$ws.Cells.Item(9, 2).Value = "FLA"
$ws.Cells.Item(5, 2).Value = "FLA"
$ws.Cells.Item(6, 2).Value = "FLA"
$ws.Cells.Item(7, 2).Value = "FLA"
$ws.Cells.Item(8, 2).Value = "FLA"

$ws.Cells.Item(9, 4).Value = "EN 1092-1 A1"
$ws.Cells.Item(5, 4).Value = "EN 1092-1 A1"
$ws.Cells.Item(6, 4).Value = "EN 1092-1 A1"
$ws.Cells.Item(7, 4).Value = "EN 1092-1 A1"
$ws.Cells.Item(8, 4).Value = "EN 1092-1 A1"

# ---------------------------------------------------------------------------
# Block 3: Blind flanges -> final rows 10..14 (PN06, PN10, PN16, PN25, PN40).
# Same PN40/06/10/16/25 typed order; column D re-uses "EN 1092-1 A1".
# ---------------------------------------------------------------------------
$ws.Cells.Item(14, 1).Value = "Blind Flange: PN40"
$ws.Cells.Item(10, 1).Value = "Blind Flange: PN06"
$ws.Cells.Item(11, 1).Value = "Blind Flange: PN10"
$ws.Cells.Item(12, 1).Value = "Blind Flange: PN16"
$ws.Cells.Item(13, 1).Value = "Blind Flange: PN25"

$ws.Cells.Item(14, 2).Value = "FLABL"
$ws.Cells.Item(10, 2).Value = "FLABL"
$ws.Cells.Item(11, 2).Value = "FLABL"
$ws.Cells.Item(12, 2).Value = "FLABL"
$ws.Cells.Item(13, 2).Value = "FLABL"

$ws.Cells.Item(14, 4).Value = "EN 1092-1 A1"
$ws.Cells.Item(10, 4).Value = "EN 1092-1 A1"
$ws.Cells.Item(11, 4).Value = "EN 1092-1 A1"
$ws.Cells.Item(12, 4).Value = "EN 1092-1 A1"
$ws.Cells.Item(13, 4).Value = "EN 1092-1 A1"

# ---------------------------------------------------------------------------
# Block 4: Elbows -> final rows 15..17 (2D, 3D, 5D). Column A typed as 3D,
# 2D, 5D; column D re-uses "EN 10253-2".
# ---------------------------------------------------------------------------
$ws.Cells.Item(16, 1).Value = "EN 10253-2 - Elbow: 3D"
$ws.Cells.Item(15, 1).Value = "EN 10253-2 - Elbow: 2D"
$ws.Cells.Item(17, 1).Value = "EN 10253-2 - Elbow: 5D"

$ws.Cells.Item(16, 2).Value = "BOG"
$ws.Cells.Item(15, 2).Value = "BOG"
$ws.Cells.Item(17, 2).Value = "BOG"

$ws.Cells.Item(16, 4).Value = "EN 10253-2"
$ws.Cells.Item(15, 4).Value = "EN 10253-2"
$ws.Cells.Item(17, 4).Value = "EN 10253-2"

# ---------------------------------------------------------------------------
# Column widths widened to fit the new, longer descriptions.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.85
$ws.Columns.Item(4).ColumnWidth = 11.7

# ---------------------------------------------------------------------------
# Final view state: scrolled down with the next empty row selected.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("E18").Select() | Out-Null
